$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45918, 1.325, 8.528),
    @(45918.01041666666, 0, 18.199),
    @(45918.02083333334, 0, 15.359),
    @(45918.03125, 0, 15.913),
    @(45918.04166666666, 0, 19.142),
    @(45918.05208333334, 0, 4.485),
    @(45918.0625, 0, 4.599),
    @(45918.07291666666, 0, 7.912),
    @(45918.08333333334, 0, 9.766),
    @(45918.09375, 0.527, 1.088),
    @(45918.10416666666, 0.105, 0.275),
    @(45918.11458333334, 0, 0.26),
    @(45918.125, 0, 1.061),
    @(45918.13541666666, 0, 6.619),
    @(45918.14583333334, 2.712, 0.548),
    @(45918.15625, 0, 5.307),
    @(45918.16666666666, 0, 10.716),
    @(45918.17708333334, 0.07000000000000001, 8.914999999999999),
    @(45918.1875, 1.035, 4.013),
    @(45918.19791666666, 0.6919999999999999, 5.424),
    @(45918.20833333334, 0, 31.023),
    @(45918.21875, 0.005, 14.424),
    @(45918.22916666666, 3.171, 1.072),
    @(45918.23958333334, 1.638, 4.658),
    @(45918.25, 0.17, 3.182),
    @(45918.26041666666, 5.116, 0.108),
    @(45918.27083333334, 1.689, 1.813),
    @(45918.28125, 3.851, 2.711),
    @(45918.29166666666, 11.82, 3.572),
    @(45918.30208333334, 0.001, 10.589),
    @(45918.3125, 0, 5.315),
    @(45918.32291666666, 0, 5.804),
    @(45918.33333333334, 0.301, 9.470000000000001),
    @(45918.34375, 0, 9.817),
    @(45918.35416666666, 0, 4.439),
    @(45918.36458333334, 0, 12.629),
    @(45918.375, 0, 37.007),
    @(45918.38541666666, 0, 23.5),
    @(45918.39583333334, 0, 9.457000000000001),
    @(45918.40625, 0, 1.739),
    @(45918.41666666666, 4.112, 0.452),
    @(45918.42708333334, 0.079, 5.33)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Apply the same style as existing A column cells (date style) to the new rows A42:A43
$ws.Range("A41").Copy() | Out-Null
$ws.Range("A42:A43").PasteSpecial(-4122) | Out-Null